$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the value for "Adez Warm Shelf Share" -> "Adez Shelf Share"
$ws.Range("A44").Value = "Adez Shelf Share"

# Delete the now-obsolete "Adez Cold Shelf Share" row entirely (row 45)
$ws.Rows("45:45").Delete()

# Update the active selection to match the target state
$ws.Range("A45").Select()
